# Update "想去人数" (wanted-to-go count) figures on the "展览" and
# "全部类型" worksheets, reflecting refreshed scrape data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 657
$ws1.Range("F4").Value  = 249
$ws1.Range("F6").Value  = 10043
$ws1.Range("F8").Value  = 905
$ws1.Range("F10").Value = 5542
$ws1.Range("F11").Value = 9
$ws1.Range("F12").Value = 10
$ws1.Range("F15").Value = 3072
$ws1.Range("F17").Value = 298
$ws1.Range("F20").Value = 9
$ws1.Range("F22").Value = 16
$ws1.Range("F23").Value = 1517

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 657
$ws4.Range("F5").Value  = 249
$ws4.Range("F7").Value  = 10043
$ws4.Range("F9").Value  = 905
$ws4.Range("F11").Value = 5542
$ws4.Range("F12").Value = 9
$ws4.Range("F13").Value = 10
$ws4.Range("F16").Value = 3072
$ws4.Range("F18").Value = 298
$ws4.Range("F21").Value = 9
$ws4.Range("F23").Value = 16
$ws4.Range("F24").Value = 1517

$wb.Save()
